# Added Test Data for UK Market
#
# Mirrors the existing "Poland" country sheet (last tab) to create a new
# "UK" tab at the end of the workbook, then tweaks the few cells that are
# country-specific: the market name in B2 and the two accessory rows that
# are listed in a different order on the Poland sheet.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Poland")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate the Poland sheet and drop it right after the last tab.
$template.Copy($null, $lastSheet)

$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# Country-specific overrides.
$uk.Range("B2").Value = "UK Market"
$uk.Range("A9").Value = "MX-DPBX"
$uk.Range("A10").Value = "MX-BBX"

# Make the new sheet the active one, matching a freshly-added tab.
$uk.Range("A10").Select()
$uk.Activate()
